$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Equipment row (row 2) -> replace old A123 equipment with new W333 equipment
$ws.Range("A2").Value = "W333"
$ws.Range("C2").Value = "W333 200017758"

# Update selection to B4
$ws.Range("B4").Select()
